# Append the 2025-10-02 profit figure as a new row at the bottom of the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use a leading apostrophe so Excel stores the date as text (matching the
# existing "MM/DD/YYYY" text entries in column A) instead of auto-converting
# it to a date serial number.
$ws.Range("A46").Value = "'10/02/2025"
$ws.Range("B46").Value = 15812.34
